$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M6").ClearContents()
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("H6").Value = 2249.8572
$ws.Range("J16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("M32").Value = -2248.5
$ws.Range("H32").Value = 2813.5715
$ws.Range("K32").Value = 2574.5
$ws.Range("I32").Value = 2574.5
$ws.Range("L48").Value = 10977
$ws.Range("J48").Value = 3659
$ws.Range("N48").Value = -11561
$ws.Range("H48").Value = 3659
$ws.Range("J51").Value = 4060.8
$ws.Range("N51").Value = -5028.8
$ws.Range("H51").Value = 2858.6667
$ws.Range("L51").Value = 4060.8
$ws.Range("N56").Value = -12045
$ws.Range("H56").Value = 3659
$ws.Range("J56").Value = 3659
$ws.Range("L56").Value = 10977
$ws.Range("H101").Value = 1778563.8
$ws.Range("J101").Value = 3000
$ws.Range("N101").Value = -12244
$ws.Range("I101").Value = 2222454.8
$ws.Range("L101").Value = 9000
$ws.Range("K101").Value = 6667364.399999999
$ws.Range("M101").Value = -6665742.399999999
$ws.Range("I107").Value = 541
$ws.Range("L107").Value = 2599
$ws.Range("M107").Value = 1379
$ws.Range("H107").Value = 1146.2941
$ws.Range("K107").Value = 541
$ws.Range("J107").Value = 2599
$ws.Range("N107").Value = -6439
$ws.Range("H129").Value = 1243.875
$ws.Range("M129").Value = 3870.5
$ws.Range("K129").Value = 1129.5
$ws.Range("I129").Value = 376.5
$ws.Range("I137").Value = 833
$ws.Range("H137").Value = 1520.3
$ws.Range("K137").Value = 2499
$ws.Range("M137").Value = 51
$ws.Range("L138").Value = 6199.5
$ws.Range("N138").Value = -16479.5
$ws.Range("J138").Value = 2066.5
$ws.Range("H138").Value = 2512.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 617886.4399999999
$ws.Range("M2").Value = -794124.9
$ws.Range("I2").Value = 794237.9
$ws.Range("K2").Value = 794237.9
$ws.Range("M4").Value = -733.5
$ws.Range("I4").Value = 849.5
$ws.Range("H4").Value = 766.5
$ws.Range("K4").Value = 849.5
$ws.Range("M6").Value = -5004827
$ws.Range("I6").Value = 5005000
$ws.Range("K6").Value = 5005000
$ws.Range("H6").Value = 5005000
$ws.Range("K14").Value = 0
$ws.Range("H14").Value = 2850
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 2850
$ws.Range("M14").ClearContents()
$ws.Range("L14").Value = 2850
$ws.Range("N14").Value = -3200
$ws.Range("L23").Value = 100007
$ws.Range("N23").Value = -100525
$ws.Range("J23").Value = 100007
$ws.Range("H23").Value = 70004.75
$ws.Range("J37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("H37").Value = 20000
$ws.Range("N37").Value = -20546
$ws.Range("I116").Value = 794237.9
$ws.Range("K116").Value = 794237.9
$ws.Range("H116").Value = 617886.4399999999
$ws.Range("M116").Value = -791943.9
$ws.Range("M122").Value = -1222.75
$ws.Range("I122").Value = 1224.25
$ws.Range("H122").Value = 2639.6
$ws.Range("K122").Value = 3672.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").Value = -794123.9
$ws.Range("H3").Value = 617886.4399999999
$ws.Range("I3").Value = 794237.9
$ws.Range("K3").Value = 794237.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3024.353
$ws.Range("I31").Value = 3097.6667
$ws.Range("K31").Value = 3097.6667
$ws.Range("M31").Value = -2802.6667
$ws.Range("K34").Value = 3097.6667
$ws.Range("I34").Value = 3097.6667
$ws.Range("H34").Value = 3024.353
$ws.Range("M34").Value = -2895.6667
$ws.Range("J141").Value = 58831.332
$ws.Range("H141").Value = 60712.57
$ws.Range("L141").Value = 58831.332
$ws.Range("N141").Value = -69191.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N131").Value = -46762.377
$ws.Range("H131").Value = 12046.8545
$ws.Range("J131").Value = 12227.459
$ws.Range("L131").Value = 36682.377

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J12").Value = 2854002.2
$ws.Range("L12").Value = 2854002.2
$ws.Range("H12").Value = 5704375.5
$ws.Range("N12").Value = -2854282.2
$ws.Range("M43").Value = -5004849
$ws.Range("N43").ClearContents()
$ws.Range("K43").Value = 5005000
$ws.Range("I43").Value = 5005000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H43").Value = 5005000
$ws.Range("J46").Value = 23950
$ws.Range("N46").Value = -24262
$ws.Range("H46").Value = 23160
$ws.Range("L46").Value = 23950
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31640
$ws.Range("J57").Value = 30000
$ws.Range("H57").Value = 30000
$ws.Range("I80").Value = 2324.25
$ws.Range("K80").Value = 2324.25
$ws.Range("H80").Value = 2511.111
$ws.Range("M80").Value = -1326.25
$ws.Range("K83").Value = 11621.25
$ws.Range("M83").Value = -6629.25
$ws.Range("H83").Value = 2511.111
$ws.Range("I83").Value = 2324.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N40").Value = -12670.5
$ws.Range("J40").Value = 12398.5
$ws.Range("M40").Value = -13183.1
$ws.Range("K40").Value = 13319.1
$ws.Range("L40").Value = 12398.5
$ws.Range("I40").Value = 13319.1
$ws.Range("H40").Value = 13056.071
$ws.Range("M122").Value = -61861.375
$ws.Range("I122").Value = 21437.125
$ws.Range("H122").Value = 21721.889
$ws.Range("K122").Value = 64311.375
$ws.Range("K132").Value = 4279.9998
$ws.Range("M132").Value = -1749.9998
$ws.Range("H132").Value = 2553.6316
$ws.Range("I132").Value = 1426.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M122").Value = -216623.11
$ws.Range("I122").Value = 73024.37
$ws.Range("H122").Value = 43273
$ws.Range("K122").Value = 219073.11
$ws.Range("H126").Value = 6691.1904
$ws.Range("I126").Value = 10119.333
$ws.Range("N126").Value = -11300.9999
$ws.Range("K126").Value = 30357.999
$ws.Range("L126").Value = 6360.999899999999
$ws.Range("J126").Value = 2120.3333
$ws.Range("M126").Value = -27887.999
$ws.Range("K132").Value = 3057
$ws.Range("M132").Value = -527
$ws.Range("H132").Value = 1451.0286
$ws.Range("I132").Value = 1019
$ws.Range("K136").Value = 64105278
$ws.Range("H136").Value = 15874273
$ws.Range("I136").Value = 21368426
$ws.Range("M136").Value = -64102728

